{"js": "// Remove the floating \"Rect\u00e1ngulo 18\" textbox (the \"{{ nombre }}\" placeholder\n// rectangle) from the document body. The paragraph that anchors the shape is\n// left in place but ends up empty once the shape is gone.\nconst body = context.document.body;\nconst shapes = body.shapes;\nshapes.load(\"items/id,items/name,items/type\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < shapes.items.length; i++) {\n  const s = shapes.items[i];\n  if (s.name === \"Rect\u00e1ngulo 18\" || s.id === 18) {\n    target = s;\n    break;\n  }\n}\n\n// Fallback: if the name/id didn't match (e.g. renumbered on re-save), find\n// the text box whose text still holds the \"{{ nombre }}\" merge placeholder.\nif (!target) {\n  const textBoxes = shapes.items.filter((s) => s.type === \"TextBox\");\n  for (const s of textBoxes) {\n    s.textFrame.textRange.load(\"text\");\n  }\n  await context.sync();\n  for (const s of textBoxes) {\n    if (s.textFrame.textRange.text.indexOf(\"nombre\") !== -1) {\n      target = s;\n      break;\n    }\n  }\n  // Last resort: just take the only text box, if there's exactly one.\n  if (!target && textBoxes.length === 1) {\n    target = textBoxes[0];\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the floating \"Rect\u00e1ngulo 18\" textbox (the \"{{ nombre }}\" placeholder\n# rectangle) from the document. The paragraph that anchors the shape stays in\n# the document but ends up empty once the shape itself is deleted.\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Shapes.Count; $i++) {\n    $s = $d.Shapes.Item($i)\n    if ($s.Name -eq \"Rect\u00e1ngulo 18\" -or $s.ID -eq 18) {\n        $target = $s\n        break\n    }\n}\n\n# Fallback: if name/id didn't match (e.g. renumbered on re-save), find the\n# text box (msoTextBox = 17) whose text still holds the merge placeholder.\nif ($target -eq $null) {\n    $textBoxes = @()\n    for ($i = 1; $i -le $d.Shapes.Count; $i++) {\n        $s = $d.Shapes.Item($i)\n        if ($s.Type -eq 17) {\n            $textBoxes += $s\n            if ($s.TextFrame.TextRange.Text -like \"*nombre*\") {\n                $target = $s\n                break\n            }\n        }\n    }\n    # Last resort: just take the only text box, if there's exactly one.\n    if ($target -eq $null -and $textBoxes.Count -eq 1) {\n        $target = $textBoxes[0]\n    }\n}\n\nif ($target -ne $null) {\n    $target.Delete()\n}\n"}
